# Update the worksheet date and the 25 division problems/answers.
$d = $word.ActiveDocument

# 1. Update the date heading paragraph.
$d.Content.Find.Execute("2023-08-28 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-08-29 Tuesday", 2)

# 2. Update each table cell by explicit (row, column) address so that
#    duplicate/overlapping old-vs-new values can't cause a Find/Replace
#    to clobber the wrong cell.
$t = $d.Tables.Item(1)

$cellUpdates = @(
    @{ Row = 1;  Col = 1; Text = "76÷7=10, 6" },
    @{ Row = 1;  Col = 2; Text = "22÷7=3, 1" },
    @{ Row = 1;  Col = 3; Text = "39÷7=5, 4" },
    @{ Row = 1;  Col = 4; Text = "12÷4=3, 0" },
    @{ Row = 1;  Col = 5; Text = "97÷8=12, 1" },

    @{ Row = 5;  Col = 1; Text = "42÷2=21, 0" },
    @{ Row = 5;  Col = 2; Text = "83÷8=10, 3" },
    @{ Row = 5;  Col = 3; Text = "67÷3=22, 1" },
    @{ Row = 5;  Col = 4; Text = "63÷3=21, 0" },
    @{ Row = 5;  Col = 5; Text = "25÷9=2, 7" },

    @{ Row = 9;  Col = 1; Text = "84÷5=16, 4" },
    @{ Row = 9;  Col = 2; Text = "63÷6=10, 3" },
    @{ Row = 9;  Col = 3; Text = "52÷7=7, 3" },
    @{ Row = 9;  Col = 4; Text = "44÷2=22, 0" },
    @{ Row = 9;  Col = 5; Text = "35÷2=17, 1" },

    @{ Row = 13; Col = 1; Text = "49÷6=8, 1" },
    @{ Row = 13; Col = 2; Text = "55÷8=6, 7" },
    @{ Row = 13; Col = 3; Text = "23÷4=5, 3" },
    @{ Row = 13; Col = 4; Text = "26÷5=5, 1" },
    @{ Row = 13; Col = 5; Text = "84÷3=28, 0" },

    @{ Row = 17; Col = 1; Text = "99÷5=19, 4" },
    @{ Row = 17; Col = 2; Text = "19÷8=2, 3" },
    @{ Row = 17; Col = 3; Text = "60÷7=8, 4" },
    @{ Row = 17; Col = 4; Text = "80÷7=11, 3" },
    @{ Row = 17; Col = 5; Text = "31÷3=10, 1" }
)

foreach ($u in $cellUpdates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
